$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 353
$newValue = 46075

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newValue
}
